$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Graphs")
$ws1.Activate()
$ws1.Range("M10").Select()

$ws2 = $wb.Worksheets.Item("CNN 1D")
$ws2.Activate()
$ws2.Range("T62").Select()

$ws4 = $wb.Worksheets.Item("MW CNN 1D")
$ws4.Activate()
$ws4.Range("L37").Select()

Write-Host "done"
